# Add two new columns (Giá thuốc / Đơn vị mặc định) with sample data to
# the "import thuoc" template sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---------------------------------------------------
$ws.Range("C1").Value = "Giá thuốc"
$ws.Range("D1").Value = "Đơn vị mặc định"

# --- Data (row 2) --------------------------------------------------------
$ws.Range("C2").Value = 100000
$ws.Range("D2").Value = "cái"

# --- Number formatting for the price column ------------------------------
# Apply to the data cell first, then the header cell, so the new cellXfs
# entries come out in the same order as the original author's workbook
# (plain+numberformat before bold+numberformat).
$ws.Range("C2").NumberFormat = "#,##0"
$ws.Range("C1").NumberFormat = "#,##0"

# --- Column widths ---------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 20.8

# --- Leave the cursor where the author last left it -----------------------
[void]$ws.Range("I9").Select()
